$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '63.266.24'
$ws.Range('D3').Value = '3.119.60'
$ws.Range('E3').Value = '  +4.32%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextCell 'D5' '586.51'
$ws.Range('E5').Value = '  +4.20%  '
Set-TextCell 'D6' '145.21'
$ws.Range('E6').Value = '  +4.70%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.113.30'
$ws.Range('E8').Value = '  +4.36%  '
$ws.Range('E9').Value = '  +2.16%  '
Set-TextCell 'D10' '0.150'
$ws.Range('E10').Value = '  +12.99%  '
Set-TextCell 'D11' '5.78'
$ws.Range('E11').Value = '  +9.32%  '
Set-TextCell 'D12' '0.467'
$ws.Range('E12').Value = '  +3.21%  '
$ws.Range('E13').Value = '  +8.16%  '
Set-TextCell 'D14' '35.67'
$ws.Range('E14').Value = '  +5.40%  '
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '3.634.96'
$ws.Range('E16').Value = '  +3.98%  '
Set-TextCell 'D17' '7.17'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '63.178.23'
$ws.Range('E18').Value = '  +6.15%  '
$ws.Range('D19').Value = '3.116.27'
$ws.Range('E19').Value = '  +4.02%  '
Set-TextCell 'D20' '467.80'
$ws.Range('E20').Value = '  +7.33%  '
Set-TextCell 'D21' '14.07'
$ws.Range('E21').Value = '  +3.70%  '
Set-TextCell 'D22' '0.726'
$ws.Range('E22').Value = '  +1.19%  '
Set-TextCell 'D23' '7.54'
$ws.Range('E23').Value = '  +6.46%  '
Set-TextCell 'D24' '13.29'
$ws.Range('E24').Value = '  -1.78%  '
Set-TextCell 'D25' '82.11'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D28' '8.35'
$ws.Range('E28').Value = '  +7.48%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D29' '2.68'
$ws.Range('E29').Value = '  +5.50%  '
$ws.Range('E30').Value = '  -0.04%  '
Set-TextCell 'D31' '6.87'
$ws.Range('E31').Value = '  +9.90%  '
Set-TextCell 'D32' '27.01'
$ws.Range('E32').Value = '  +4.74%  '
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').Value = '0.0₃0873'
$ws.Range('E34').Value = '  +12.15%  '
$ws.Range('E35').Value = '  +16.49%  '
$ws.Range('E36').Value = '  +5.10%  '
Set-TextCell 'D37' '3.33'
$ws.Range('E37').Value = '  +19.53%  '
Set-TextCell 'D38' '6.02'
$ws.Range('E38').Value = '  +2.31%  '
Set-TextCell 'D39' '50.89'
$ws.Range('E39').Value = '  +3.80%  '
Set-TextCell 'D40' '432.72'
$ws.Range('E40').Value = '  +8.18%  '
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').Value = '2.927.93'
$ws.Range('E42').Value = '  +5.91%  '
$ws.Range('E43').Value = '  +4.67%  '
Set-TextCell 'D44' '0.280'
$ws.Range('E44').Value = '  +11.34%  '
Set-TextCell 'D45' '0.111'
$ws.Range('E45').Value = '  +5.56%  '
$ws.Range('E46').Value = '  +7.65%  '
Set-TextCell 'D47' '35.39'
$ws.Range('E47').Value = '  +2.81%  '
Set-TextCell 'D49' '123.23'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('E51').Value = '  +4.20%  '
